# Center images in book-type files: apply the "Compact" paragraph style
# to every paragraph that consists solely of an inline picture (drawing),
# mirroring the author's markdown->docx generator change.
$d = $word.ActiveDocument

for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
    $shape = $d.InlineShapes.Item($i)
    $para = $shape.Range.Paragraphs.Item(1)
    $para.Style = "Compact"
}
